$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Title timestamp text change
# ---------------------------------------------------------------
$d.Content.Find.Execute("2025-03-01 17:37:09.768835", $false, $false, $false, $false, $false, $true, 1, $false, "2025-03-01 17:43:55.105921", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Heading style changes: Heading2/3/4/5 -> Heading1
# ---------------------------------------------------------------
$d.Paragraphs(3).Style = "Heading1"   # 08:55:00           (was Heading2)
$d.Paragraphs(4).Style = "Heading1"   # RMF-Client01        (was Heading3)
$d.Paragraphs(5).Style = "Heading1"   # Debian 6.1.128-1    (was Heading4)
$d.Paragraphs(6).Style = "Heading1"   # 10.0.0.20           (was Heading5)

# ---------------------------------------------------------------
# 3. Global whitespace tweaks: 4 leading spaces -> 2 leading spaces
#    for every "Current version:" / "Update version:" line.
# ---------------------------------------------------------------
$d.Content.Find.Execute("    Current version: ", $false, $false, $false, $false, $false, $true, 1, $false, "  Current version: ", 2) | Out-Null
$d.Content.Find.Execute("    Update version: ", $false, $false, $false, $false, $false, $true, 1, $false, "  Update version: ", 2) | Out-Null

# ---------------------------------------------------------------
# 4. Structural edits - process from the bottom of the document
#    upward so earlier (still-unprocessed) paragraph indices stay
#    valid while later ones are split/merged/removed.
# ---------------------------------------------------------------

$vtab = [char]11   # manual line break (<w:br/>) inside Word's Range.Text model

# --- Section 12 (openssh-client) - last paragraph in the document ---
# Paragraph 65: "    Affiliated CVES: None" (no break, end of doc) -> remove entirely
# Paragraph 64: "  Update version: 1:9.2p1-2+deb12u5 " + break -> becomes
#               "  Update version: 1:9.2p1-2+deb12u5" with no trailing space/break
$d.Paragraphs(65).Range.Delete() | Out-Null
$p64 = $d.Paragraphs(64)
$r64 = $p64.Range
$r64trim = $d.Range($r64.Start, $r64.End - 1)
$r64trim.Text = "  Update version: 1:9.2p1-2+deb12u5"

# --- Section 11,10,9 CVE split (paragraph 60,55,50) ---
foreach ($idx in 60,55,50) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $d.Paragraphs($idx).Range.Find.Execute("    Affiliated CVES: CVE-2013-3667, CVE-2025-1228 ", $false, $false, $false, $false, $false, $true, 1, $false, "  Affiliated CVES: ", 2) | Out-Null
    $d.Paragraphs($idx + 1).Range.InsertBefore("   CVE-2013-3667" + $vtab)
}

# --- Section 8 (libtasn1-6) ---
# Paragraph 45: "    Affiliated CVES: None " (with break) -> remove entirely
# Paragraph 44: "  Update version: 4.19.0-2+deb12u1 " + break -> becomes
#               "  Update version: Not mentioned " + break
$d.Paragraphs(45).Range.Delete() | Out-Null
$d.Paragraphs(44).Range.Find.Execute("  Update version: 4.19.0-2+deb12u1 ", $false, $false, $false, $false, $false, $true, 1, $false, "  Update version: Not mentioned ", 2) | Out-Null

# --- Sections 7,6,5,4,3,2,1 CVE split (paragraph 40,35,30,25,20,15,10) ---
foreach ($idx in 40,35,30,25,20,15,10) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $d.Paragraphs($idx).Range.Find.Execute("    Affiliated CVES: CVE-2013-3667, CVE-2025-1228 ", $false, $false, $false, $false, $false, $true, 1, $false, "  Affiliated CVES: ", 2) | Out-Null
    $d.Paragraphs($idx + 1).Range.InsertBefore("   CVE-2013-3667" + $vtab)
}

Write-Output "done"
